$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2012.215366666667
$ws.Range("H2").Value = 6036.6461
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 25.94532566666667
$ws.Range("N2").Value = 77.835977
$ws.Range("O2").Value = 0.5401813355606462
$ws.Range("P2").Value = 0.5401813355606462
$ws.Range("Q2").Value = 52207.58299963774
$ws.Range("R2").Value = 469868.2469967397
$ws.Range("S2").Value = 0.5401813355606462
$ws.Range("T2").Value = 0.5401813355606462

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2012.215366666667
$ws.Range("H3").Value = 6036.6461
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 18.51427066666667
$ws.Range("N3").Value = 55.542812
$ws.Range("O3").Value = 0.3854668692210787
$ws.Range("P3").Value = 0.3854668692210786
$ws.Range("Q3").Value = 37254.69993809258
$ws.Range("R3").Value = 335292.2994428332
$ws.Range("S3").Value = 0.3854668692210787
$ws.Range("T3").Value = 0.3854668692210786

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2012.215366666667
$ws.Range("H4").Value = 6036.6461
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.571174000000001
$ws.Range("N4").Value = 10.713522
$ws.Range("O4").Value = 0.07435179521827505
$ws.Range("P4").Value = 0.07435179521827504
$ws.Range("Q4").Value = 7185.971199840467
$ws.Range("R4").Value = 64673.74079856421
$ws.Range("S4").Value = 0.07435179521827505
$ws.Range("T4").Value = 0.07435179521827504
